$d = $word.ActiveDocument

# --- 1. Split the first paragraph's run and append the red "(This is a
#        change - Version for branch alternate)" text. ---

$p1 = $d.Paragraphs(1)
$p1Range = $p1.Range
# Trim the paragraph mark off the end of the range so we only touch the text.
$p1TextRange = $d.Range($p1Range.Start, $p1Range.End - 1)
$p1TextRange.Text = "This is a Microsoft word document.  "

$redColor = 192  # OLE color 0x0000C0 == RGB(192,0,0) == OOXML w:color C00000

$insertPoint = $d.Range($p1TextRange.End, $p1TextRange.End)
$r2Text = "(This is a change " + [char]0x2013 + " Ve"
$insertPoint.InsertAfter($r2Text)
$r2Range = $d.Range($p1TextRange.End, $p1TextRange.End + $r2Text.Length)
$r2Range.Font.Color = $redColor

$afterR2 = $r2Range.End
$r3Text = "rsion for branch alternate"
$insertPoint2 = $d.Range($afterR2, $afterR2)
$insertPoint2.InsertAfter($r3Text)
$r3Range = $d.Range($afterR2, $afterR2 + $r3Text.Length)
$r3Range.Font.Color = $redColor

$afterR3 = $r3Range.End
$r4Text = ")"
$insertPoint3 = $d.Range($afterR3, $afterR3)
$insertPoint3.InsertAfter($r4Text)
$r4Range = $d.Range($afterR3, $afterR3 + $r4Text.Length)
$r4Range.Font.Color = $redColor

# --- 2. Add a new empty shaded paragraph at the very end of the document,
#        right after the "Free at last..." paragraph. ---

$endRange = $d.Content
$endRange.Collapse(0)
$null = $endRange.InsertParagraphAfter()

$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$newParaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
  '<w:pPr><w:shd w:val="clear" w:color="auto" w:fill="F9F9F9"/></w:pPr></w:p>'
$null = $lastPara.Range.InsertXML($newParaXml)
